$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style + row height) of the last existing data row (181) down to the new rows (182:221)
$ws.Range("A181:C181").Copy() | Out-Null
$ws.Range("A182:C221").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the new tombola draw data (rows 182:221)
$data = @(
  @(45873,1,47),
  @(45873,2,85),
  @(45873,3,52),
  @(45873,4,8),
  @(45873,5,69),
  @(45873,6,99),
  @(45873,7,46),
  @(45873,8,4),
  @(45873,9,81),
  @(45873,10,16),
  @(45874,1,5),
  @(45874,2,50),
  @(45874,3,12),
  @(45874,4,37),
  @(45874,5,64),
  @(45874,6,18),
  @(45874,7,82),
  @(45874,8,41),
  @(45874,9,48),
  @(45874,10,6),
  @(45875,1,24),
  @(45875,2,24),
  @(45875,3,97),
  @(45875,4,59),
  @(45875,5,76),
  @(45875,6,19),
  @(45875,7,10),
  @(45875,8,16),
  @(45875,9,1),
  @(45875,10,57),
  @(45876,1,66),
  @(45876,2,9),
  @(45876,3,76),
  @(45876,4,39),
  @(45876,5,52),
  @(45876,6,32),
  @(45876,7,9),
  @(45876,8,49),
  @(45876,9,33),
  @(45876,10,99)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = 182 + $i
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
  $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Ensure correct row height matches the rest of the sheet (15.75, custom height)
$ws.Range("A182:C221").RowHeight = 15.75

# Update selection to match the new active cell after data entry (C222), mirroring the saved view state
$ws.Range("C222").Select()
